$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 97224760
$ws.Range("I43").Value = 200001890
$ws.Range("J43").Value = 23812524
$ws.Range("K43").Value = 200001890
$ws.Range("L43").Value = 23812524
$ws.Range("M43").Value = -200001821
$ws.Range("N43").Value = -23812662

$ws.Range("H112").Value = 1698.65
$ws.Range("J112").Value = 1986.125
$ws.Range("L112").Value = 5958.375
$ws.Range("N112").Value = -8174.375

$ws.Range("H116").Value = 2195.3
$ws.Range("J116").Value = 2438.6155
$ws.Range("L116").Value = 2438.6155
$ws.Range("N116").Value = -9322.6155

$ws.Range("H121").Value = 1955
$ws.Range("J121").Value = 1955
$ws.Range("L121").Value = 5865
$ws.Range("N121").Value = -9359

$ws.Range("H135").Value = 1369.9474
$ws.Range("I135").Value = 1119.3529
$ws.Range("J135").Value = 3500
$ws.Range("K135").Value = 10074.1761
$ws.Range("L135").Value = 31500
$ws.Range("M135").Value = -7539.176100000001
$ws.Range("N135").Value = -36570

$ws.Range("H137").Value = 1953
$ws.Range("I137").Value = 1800.4667
$ws.Range("J137").Value = 2525
$ws.Range("K137").Value = 5401.4001
$ws.Range("L137").Value = 7575
$ws.Range("M137").Value = -2851.4001
$ws.Range("N137").Value = -12675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20229.436
$ws.Range("I32").Value = 22694.725
$ws.Range("J32").Value = 8799.454
$ws.Range("K32").Value = 22694.725
$ws.Range("L32").Value = 8799.454
$ws.Range("M32").Value = -22407.725
$ws.Range("N32").Value = -9373.454

$ws.Range("H45").Value = 905.8823
$ws.Range("I45").Value = 960
$ws.Range("J45").Value = 828.5714
$ws.Range("K45").Value = 960
$ws.Range("L45").Value = 828.5714
$ws.Range("M45").Value = -583
$ws.Range("N45").Value = -1582.5714

$ws.Range("H61").Value = 2260.875
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 3241.5557
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 3241.5557
$ws.Range("M61").Value = -788
$ws.Range("N61").Value = -3665.5557

$ws.Range("H74").Value = 1168.5714
$ws.Range("I74").Value = 1013.3333
$ws.Range("J74").Value = 2100
$ws.Range("K74").Value = 1013.3333
$ws.Range("L74").Value = 2100
$ws.Range("M74").Value = -139.3333
$ws.Range("N74").Value = -3848

$ws.Range("H77").Value = 1168.5714
$ws.Range("I77").Value = 1013.3333
$ws.Range("J77").Value = 2100
$ws.Range("K77").Value = 5066.6665
$ws.Range("L77").Value = 10500
$ws.Range("M77").Value = -698.6665000000003
$ws.Range("N77").Value = -19236

$ws.Range("H136").Value = 2260.875
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 3241.5557
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 9724.667099999999
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -14824.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6253392.5
$ws.Range("I31").Value = 2363.9312
$ws.Range("J31").Value = 66680000
$ws.Range("K31").Value = 2363.9312
$ws.Range("L31").Value = 66680000
$ws.Range("M31").Value = -2068.9312
$ws.Range("N31").Value = -66680590

$ws.Range("H34").Value = 6253392.5
$ws.Range("I34").Value = 2363.9312
$ws.Range("J34").Value = 66680000
$ws.Range("K34").Value = 2363.9312
$ws.Range("L34").Value = 66680000
$ws.Range("M34").Value = -2161.9312
$ws.Range("N34").Value = -66680404

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 977.1429000000001
$ws.Range("I5").Value = 515
$ws.Range("J5").Value = 1323.75
$ws.Range("K5").Value = 1545
$ws.Range("L5").Value = 3971.25
$ws.Range("M5").Value = -1433
$ws.Range("N5").Value = -4195.25

$ws.Range("H9").Value = 12826794
$ws.Range("I9").Value = 4000
$ws.Range("J9").Value = 13339706
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 40019118
$ws.Range("M9").Value = -11776
$ws.Range("N9").Value = -40019566

$ws.Range("H32").Value = 1466.6666
$ws.Range("J32").Value = 1466.6666
$ws.Range("L32").Value = 4399.9998
$ws.Range("N32").Value = -4965.9998

$ws.Range("H122").Value = 1508.7931
$ws.Range("J122").Value = 1398.2632
$ws.Range("L122").Value = 12584.3688
$ws.Range("N122").Value = -17484.3688

$ws.Range("H131").Value = 3588111.2
$ws.Range("J131").Value = 6173689
$ws.Range("L131").Value = 18521067
$ws.Range("N131").Value = -18531147

$ws.Range("H135").Value = 977.1429000000001
$ws.Range("I135").Value = 515
$ws.Range("J135").Value = 1323.75
$ws.Range("K135").Value = 4635
$ws.Range("L135").Value = 11913.75
$ws.Range("M135").Value = -2100
$ws.Range("N135").Value = -16983.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1260.1471
$ws.Range("I102").Value = 1159.52
$ws.Range("J102").Value = 1539.6666
$ws.Range("K102").Value = 1159.52
$ws.Range("L102").Value = 1539.6666
$ws.Range("M102").Value = 462.48
$ws.Range("N102").Value = -4783.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1684.3939
$ws.Range("I46").Value = 1457.2142
$ws.Range("J46").Value = 1851.7894
$ws.Range("K46").Value = 1457.2142
$ws.Range("L46").Value = 1851.7894
$ws.Range("M46").Value = -1269.2142
$ws.Range("N46").Value = -2227.7894

$ws.Range("H132").Value = 26785.45
$ws.Range("I132").Value = 34320.535
$ws.Range("K132").Value = 102961.605
$ws.Range("M132").Value = -100431.605

$ws.Range("H136").Value = 6453
$ws.Range("I136").Value = 8075.643
$ws.Range("J136").Value = 2666.8333
$ws.Range("K136").Value = 24226.929
$ws.Range("L136").Value = 8000.499899999999
$ws.Range("M136").Value = -21676.929
$ws.Range("N136").Value = -13100.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2899.5715
$ws.Range("I132").Value = 2059.8
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6179.400000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3649.400000000001
$ws.Range("N132").Value = -20057

$ws.Range("H135").Value = 58826.5
$ws.Range("J135").Value = 58826.5
$ws.Range("L135").Value = 58826.5
$ws.Range("N135").Value = -68966.5

$ws.Range("H136").Value = 10124.143
$ws.Range("I136").Value = 12267.091
$ws.Range("J136").Value = 2266.6667
$ws.Range("K136").Value = 36801.273
$ws.Range("L136").Value = 6800.000100000001
$ws.Range("M136").Value = -34251.273
$ws.Range("N136").Value = -11900.0001

$ws.Range("H139").Value = 46759
$ws.Range("J139").Value = 48286.25
$ws.Range("L139").Value = 48286.25
$ws.Range("N139").Value = -58566.25
